$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a pure number by
# Excels type inference (losing formatting like trailing zeros, e.g. "1.00" -> 1).
# Pre-format them as Text so the literal string is preserved, matching the source data.
$textCells = @("D5","D6","D8","D9","D11","D14","D15","D18","D19","D20","D21","D23","D24","D28","D29","D30","D31","D32","D33","D34","D36","D37","D38","D40","D41","D42","D43","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (prices, volume %, and the swapped FirstDigitalUSD/Mantle rows).
$ws.Range("D2").Value = "57.193.95"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "3.063.45"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "513.25"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "142.40"
$ws.Range("E6").Value = "  +7.63%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.435"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").Value = "  +5.92%  "
$ws.Range("D12").Value = "3.599.14"
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "25.75"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "0.0000164"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "57.393.81"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "3.066.55"
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").Value = "6.09"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "8.15"
$ws.Range("E20").Value = "  +6.22%  "
$ws.Range("D21").Value = "335.28"
$ws.Range("E21").Value = "  +6.96%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "0.498"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").Value = "65.32"
$ws.Range("E24").Value = "  +4.14%  "
$ws.Range("E25").Value = "  +7.16%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "0.0₃0930"
$ws.Range("E27").Value = "  +11.46%  "
$ws.Range("D28").Value = "6.41"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").Value = "7.08"
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").Value = "1.80"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").Value = "20.69"
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "154.51"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").Value = "26.12"
$ws.Range("E36").Value = "  +9.06%  "
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("D38").Value = "0.0680"
$ws.Range("E38").Value = "  +5.69%  "
$ws.Range("D39").Value = "3.107.53"
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").Value = "36.70"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.667"
$ws.Range("E43").Value = "  +4.82%  "
$ws.Range("D44").Value = "2.266.11"
$ws.Range("E44").Value = "  +7.55%  "
$ws.Range("E45").Value = "  +8.06%  "
$ws.Range("D46").Value = "1.38"
$ws.Range("E46").Value = "  +4.39%  "
$ws.Range("D47").Value = "0.951"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("D48").Value = "20.31"
$ws.Range("E48").Value = "  +9.54%  "
$ws.Range("D49").Value = "5.85"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("D51").Value = "0.689"
$ws.Range("E51").Value = "  +7.41%  "
